# Saisie des premieres mesures pour l'arbre piege de la ligne 11 (debut BMP)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = 38
$ws.Range("E11").Value = 36
$ws.Range("H11").Value = 15.1

# La derniere cellule active/selectionnee passe de D11 a F11
$ws.Range("F11").Select()
